# VH-00 - Ajustes do Back-end
#
# Source sheet ("Exportação") stays structurally the same; only the
# selected cell / scroll position changed.
#
# Target sheet ("Importação") gains a new "CNPJ" field size/position
# (row 4) and a whole new copied block (rows 33-40) that mirrors the
# header + product-field rows (13-20), with the last three fields
# (CATEGORIA / SUBCATEGORIA / DESCRIÇÃO) re-sized.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Exportação
$ws2 = $wb.Worksheets.Item(2)   # Importação

# ---------------------------------------------------------------
# Importação: field 3 (CNPJ) Tamanho/Posição change
#   Tamanho  "18"      -> "14"
#   Posição  "033-050" -> "033-046"
# ---------------------------------------------------------------
$ws2.Range("C4").Value = "14"
$ws2.Range("D4").Value = "033-046"

# ---------------------------------------------------------------
# Importação: duplicate the header + 7 data rows (A13:F20) down to
# A33:F40, carrying styles/number formats along, then touch up the
# three rows whose Tamanho/Posição differ from the source block.
# ---------------------------------------------------------------
$ws2.Range("A13:F20").Copy($ws2.Range("A33"))

$ws2.Range("C38").Value = "015"
$ws2.Range("D38").Value = "043-057"

$ws2.Range("C39").Value = "015"
$ws2.Range("D39").Value = "058-072"

$ws2.Range("C40").Value = "85"
$ws2.Range("D40").Value = "073-157"

# ---------------------------------------------------------------
# View state: Exportação scrolled with C10 selected; Importação
# scrolled further down with A33:F40 selected (F40 active) and stays
# the active tab.
# ---------------------------------------------------------------
$ws1.Activate()
$ws1.Range("C10").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1

$ws2.Activate()
$ws2.Range("A33:F40").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
